$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '36.294.45'
$ws.Range("E2").Value = '  -1.52%  '

# Row 3
$ws.Range("D3").Value = '2.032.40'
$ws.Range("E3").Value = '  -0.88%  '

# Row 4
$ws.Range("E4").Value = '  +0.01%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '243.87'
$ws.Range("E5").Value = '  -0.50%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.656'
$ws.Range("E6").Value = '  +0.31%  '

# Row 7
$ws.Range("E7").Value = '  +0.04%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '52.73'
$ws.Range("E8").Value = '  -8.64%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '61.50'
$ws.Range("E9").Value = '  +5.15%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.358'
$ws.Range("E10").Value = '  -2.53%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0736'
$ws.Range("E11").Value = '  -4.89%  '

# Row 12
$ws.Range("E12").Value = '  -4.11%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.939'
$ws.Range("E13").Value = '  +8.01%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '14.23'
$ws.Range("E14").Value = '  -6.13%  '

# Row 15
$ws.Range("D15").Value = '2.329.70'
$ws.Range("E15").Value = '  -0.75%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '5.27'
$ws.Range("E16").Value = '  -5.53%  '

# Row 17
$ws.Range("D17").Value = '2.040.69'
$ws.Range("E17").Value = '  -0.16%  '

# Row 18
$ws.Range("D18").Value = '36.038.41'
$ws.Range("E18").Value = '  -2.08%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '16.62'
$ws.Range("E19").Value = '  -6.23%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '70.73'
$ws.Range("E20").Value = '  -3.42%  '

# Row 21
$ws.Range("D21").Value = '0.0₃0842'
$ws.Range("E21").Value = '  -4.74%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '235.98'
$ws.Range("E22").Value = '  +0.06%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.11'
$ws.Range("E23").Value = '  -4.79%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.999'
$ws.Range("E24").Value = '  -0.24%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.36'
$ws.Range("E25").Value = '  -3.67%  '

# Row 26
$ws.Range("E26").Value = '  -0.67%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '162.62'
$ws.Range("E27").Value = '  -3.46%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.04'
$ws.Range("E28").Value = '  -13.01%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '19.67'
$ws.Range("E29").Value = '  -1.25%  '

# Row 30
$ws.Range("E30").Value = '  -3.55%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.15'
$ws.Range("E31").Value = '  +4.27%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.86'
$ws.Range("E32").Value = '  -10.37%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0587'
$ws.Range("E33").Value = '  -3.89%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.33'
$ws.Range("E34").Value = '  -10.37%  '

# Row 35
$ws.Range("E35").Value = '  +0.09%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0847'
$ws.Range("E36").Value = '  +7.25%  '

# Row 37
$ws.Range("E37").Value = '  -1.30%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.18'
$ws.Range("E38").Value = '  -6.28%  '

# Row 39
$ws.Range("B39").Value = 'TrustWalletToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.22'
$ws.Range("E39").Value = '  -7.18%  '

# Row 40
$ws.Range("B40").Value = 'THORChain'
$ws.Range("C40").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '4.83'
$ws.Range("E40").Value = '  -3.23%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.85'
$ws.Range("E41").Value = '  -4.94%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0211'
$ws.Range("E42").Value = '  -4.98%  '

# Row 43
$ws.Range("E43").Value = '  -5.15%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '92.05'
$ws.Range("E44").Value = '  -4.22%  '

# Row 45
$ws.Range("B45").Value = 'Maker'
$ws.Range("C45").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D45").Value = '1.380.95'
$ws.Range("E45").Value = '  +5.48%  '

# Row 46
$ws.Range("B46").Value = 'Cronos'
$ws.Range("C46").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0885'
$ws.Range("E46").Value = '  -5.71%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '7.42'
$ws.Range("E47").Value = '  +10.44%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '15.60'
$ws.Range("E48").Value = '  -7.32%  '

# Row 49
$ws.Range("E49").Value = '  +2.08%  '

# Row 50
$ws.Range("D50").Value = '2.219.36'
$ws.Range("E50").Value = '  -0.58%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.23'
$ws.Range("E51").Value = '  -5.24%  '
